# Insert a new data row at row 93 (shifts existing rows 93..230 down to 94..231)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(93).Insert()

$ws.Cells.Item(93, 1).Value = 10
$ws.Cells.Item(93, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(93, 3).Value = "La Araucanía"
$ws.Cells.Item(93, 4).Value = 44792
$ws.Cells.Item(93, 5).Value = 9
$ws.Cells.Item(93, 6).Value = 100112005
$ws.Cells.Item(93, 7).Value = "Puerro"
$ws.Cells.Item(93, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 50
$ws.Cells.Item(93, 11).Value = 18000
$ws.Cells.Item(93, 12).Value = 18000
$ws.Cells.Item(93, 13).Value = 18000
$ws.Cells.Item(93, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(93, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(93, 16).Value = 1500
$ws.Cells.Item(93, 17).Value = 12
$ws.Cells.Item(93, 18).Value = "Hortaliza"
